# Applies the "add model to tg" edit:
#  - B2 gains a trailing space
#  - B3 is reduced to a single leading/trailing space (no embedded newline)
#  - B4/C4/D4/E4 are trimmed down to only the "СЛУЧАЙ" fragment (the rest of
#    the original multi-paragraph text in row 4 is split out into two new
#    rows: row 5 ("Группа создана ...") and row 6 ("Как настоящий ...")
#  - Two new rows (5 and 6) are appended, copying the formatting of column A
#    from row 4 (bold / centered / bordered style), with F set to 0.
#  - The sheet dimension grows from A1:F4 to A1:F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 -----------------------------------------------------------
$ws.Range("B2").Value = "Нет, хуйня это ваш зумерский снюс. "

# ---- Row 3 -----------------------------------------------------------
$ws.Range("B3").Value = " Всем привет, как дела? "

# ---- Row 4 (trimmed to just the ЧАЙ fragment) -------------------------
$ws.Range("B4").Value = "  ║ [ Н҉А҉ С҉Л҉У҉Ч҉А҉Й҉ Е҉С҉Л҉И҉ Я҉ У҉М҉Р҉У.  "
$ws.Range("C4").Value = "║ н҉а҉ с҉л҉у҉ч҉а҉й҉ е҉с҉л҉и҉ я҉ у҉м҉р҉у"
$ws.Range("D4").Value = "║ н҉а҉ с҉л҉у҉ч҉а҉й҉ е҉с҉л҉и҉ я҉ у҉м҉р҉"
$ws.Range("E4").Value = "║ н҉а҉ с҉л҉у҉ч҉а҉й҉ е҉с҉л҉и҉ я҉ у҉м҉р҉"
$ws.Range("F4").Value = 0

# ---- Row 5 (new) -------------------------------------------------------
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = " Группа создана в экспериментальных целях "
$ws.Range("C5").Value = "группа создана в экспериментальных целях"
$ws.Range("D5").Value = "группа создать в экспериментальный цель"
$ws.Range("E5").Value = "группа создать экспериментальный цель"
$ws.Range("F5").Value = 0

# ---- Row 6 (new) -------------------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = " Как настоящий ру`$`$кий я предпочту убойный насвай, синтетический гашиш, метадон, и мефедрон"
$ws.Range("C6").Value = "как настоящий ру кий я предпочту убойный насвай синтетический гашиш метадон и мефедрон"
$ws.Range("D6").Value = "как настоящий ру кий я предпочесть убойный насвая синтетический гашиш метадон и мефедрон"
$ws.Range("E6").Value = "настоящий ру кий предпочесть убойный насвая синтетический гашиш метадон мефедрон"
$ws.Range("F6").Value = 0

# Copy column-A formatting (bold / centered / thin border) from row 4 down
# to the two freshly created rows, matching the style used on A2:A4.
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "edit complete"
